$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: for Price column values that look numeric (single decimal
# point, e.g. "241.14"), Excel auto-converts Range.Value assignments to a
# Number. The source data keeps these as plain text (matching the rest of the
# Price column, including multi-dot values like "29.157.50" which can only be
# text). Forcing NumberFormat to text ("@") before the assignment keeps the
# cell a Text cell; ClearFormats() afterwards drops the temporary number format
# so the cell style reverts to the original default (no explicit style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Update Price (D) and Volume(1h) (E) columns for rows with refreshed crypto data
$ws.Range("D2").Value = "29.157.50"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.823.27"
$ws.Range("E3").Value = "  -0.46%  "
Set-TextValue $ws.Range("D4") "0.9989"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "241.14"
$ws.Range("E5").Value = "  -0.76%  "
Set-TextValue $ws.Range("D6") "0.6192"
$ws.Range("E6").Value = "  -1.40%  "
Set-TextValue $ws.Range("D7") "0.9995"
$ws.Range("E7").Value = "  -0.08%  "
Set-TextValue $ws.Range("D8") "0.07337"
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("E9").Value = "  -1.04%  "
Set-TextValue $ws.Range("D10") "22.98"
$ws.Range("E10").Value = "  -1.04%  "
Set-TextValue $ws.Range("D11") "0.07662"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.824.24"
$ws.Range("E12").Value = "  -0.45%  "
Set-TextValue $ws.Range("D13") "4.958"
$ws.Range("E13").Value = "  -1.36%  "
Set-TextValue $ws.Range("D14") "0.6612"
$ws.Range("E14").Value = "  -1.19%  "
Set-TextValue $ws.Range("D15") "82.13"
$ws.Range("E15").Value = "  -0.89%  "
Set-TextValue $ws.Range("D16") "0.000008921"
$ws.Range("E16").Value = "  -4.92%  "
Set-TextValue $ws.Range("D17") "5.829"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "29.138.92"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "2.070.69"
$ws.Range("E19").Value = "  -0.40%  "
Set-TextValue $ws.Range("D20") "236.91"
$ws.Range("E20").Value = "  +6.15%  "
Set-TextValue $ws.Range("D21") "12.40"
$ws.Range("E21").Value = "  -1.49%  "
Set-TextValue $ws.Range("D22") "0.9996"
$ws.Range("E22").Value = "  -0.28%  "
Set-TextValue $ws.Range("D23") "7.185"
$ws.Range("E23").Value = "  +0.56%  "
Set-TextValue $ws.Range("D24") "1.000"
$ws.Range("E24").Value = "  +0.01%  "
Set-TextValue $ws.Range("D25") "158.20"
$ws.Range("E25").Value = "  -1.18%  "
Set-TextValue $ws.Range("D26") "0.1416"
$ws.Range("E26").Value = "  +1.33%  "
Set-TextValue $ws.Range("D27") "8.431"
$ws.Range("E27").Value = "  -0.83%  "
Set-TextValue $ws.Range("D28") "17.63"
$ws.Range("E28").Value = "  -1.43%  "
Set-TextValue $ws.Range("D29") "1.481"
$ws.Range("E29").Value = "  -0.98%  "
Set-TextValue $ws.Range("D30") "0.05564"
$ws.Range("E30").Value = "  -4.51%  "
Set-TextValue $ws.Range("D31") "4.087"
$ws.Range("E31").Value = "  -0.78%  "
Set-TextValue $ws.Range("D32") "4.089"
$ws.Range("E32").Value = "  -1.76%  "
Set-TextValue $ws.Range("D33") "1.202"
$ws.Range("E33").Value = "  -0.07%  "
Set-TextValue $ws.Range("D34") "1.829"
$ws.Range("E34").Value = "  -0.26%  "
Set-TextValue $ws.Range("D35") "0.7329"
$ws.Range("E35").Value = "  -1.13%  "
Set-TextValue $ws.Range("D36") "1.129"
$ws.Range("E36").Value = "  -0.87%  "
Set-TextValue $ws.Range("D37") "2.619"
$ws.Range("E37").Value = "  -1.81%  "
Set-TextValue $ws.Range("D38") "2.834"
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("D39").Value = "1.213.62"
$ws.Range("E39").Value = "  -0.96%  "
Set-TextValue $ws.Range("D40") "0.01758"
$ws.Range("E40").Value = "  -1.28%  "

# Row 41/42: FraxShare and TrustWalletToken swapped positions, with updated data
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "0.9245"
$ws.Range("E41").Value = "  +3.67%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "6.302"
$ws.Range("E42").Value = "  -2.97%  "

Set-TextValue $ws.Range("D43") "0.9991"
$ws.Range("E43").Value = "  -0.11%  "
Set-TextValue $ws.Range("D44") "101.37"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "1.969.90"
$ws.Range("E45").Value = "  -0.33%  "
Set-TextValue $ws.Range("D46") "64.68"
$ws.Range("E46").Value = "  -1.87%  "
Set-TextValue $ws.Range("D47") "0.5078"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  -2.72%  "
Set-TextValue $ws.Range("D49") "0.4004"
$ws.Range("E49").Value = "  -1.57%  "
Set-TextValue $ws.Range("D50") "9.065"
$ws.Range("E50").Value = "  +0.77%  "
Set-TextValue $ws.Range("D51") "0.05754"
$ws.Range("E51").Value = "  -1.13%  "
